$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates:
#   D2: "Study" -> "study"
#   G2: "Original_Space" -> "CoordSys"
$ws.Range("D2").Value = "study"
$ws.Range("G2").Value = "CoordSys"

# Data rows for Spengler_2009 / Nahab_in_press (rows 34-45) used the
# Talairach coordinate system, labeled "TAL"; relabel as "T88".
$ws.Range("G34:G45").Value = "T88"

# Match the author's last active-cell selection.
$ws.Range("G25").Select() | Out-Null
